$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2096
$ws1.Range("F9").Value = 10800
$ws1.Range("F15").Value = 9015
$ws1.Range("F16").Value = 1119
$ws1.Range("F18").Value = 5293
$ws1.Range("F20").Value = 3358

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2096
$ws4.Range("F12").Value = 10800
$ws4.Range("F18").Value = 9015
$ws4.Range("F19").Value = 1119
$ws4.Range("F21").Value = 5293
$ws4.Range("F23").Value = 3358
